$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5 - new test case "roundTripFlightSearch" (Dhaka -> Chittagong)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "roundTripFlightSearch"
$ws.Range("B5").Value = "Dhaka"
$ws.Range("C5").Value = "Chittagong"

$ws.Range("D5").Value = "22 December 2022"
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E5").Value = "'28 December 2022"

$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = "Economy"
$ws.Range("J5").Value = "Round-trip"

# K5 should not exist on this row any more
$ws.Range("K5").Clear()

# ---------------------------------------------------------------------------
# Row 6 - continuation of "roundTripFlightSearch" (Dhaka -> Kolkata)
# ---------------------------------------------------------------------------
# A6 should not exist on this row any more
$ws.Range("A6").Clear()

$ws.Range("B6").Value = "Dhaka"

$ws.Range("D6").Value = "22 January 2023"
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E6").Value = "'28 March 2023"

$ws.Range("C6").Value = "Kolkata"

$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = "Economy"
$ws.Range("J6").Value = "Round-trip"
$ws.Range("K6").Value = "roundTripFlightSearch"

# ---------------------------------------------------------------------------
# Row 7 - new test case (Toronto -> Calgary)
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Toronto"
$ws.Range("C7").Value = "Calgary"

$ws.Range("D7").Value = "'12 January 2023"
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("E7").Value = "'13 January 2023"
$ws.Range("D2").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Economy"
$ws.Range("J7").Value = "Round-trip"

# K7 should not exist on this row any more
$ws.Range("K7").Clear()

# ---------------------------------------------------------------------------
# Row 8 - K8 no longer present
# ---------------------------------------------------------------------------
$ws.Range("K8").Clear()

# Leave the selection where the last edit happened
$ws.Range("K6").Select()
